$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 140; all existing rows from 140 down shift to 141+.
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with the new weekly price record
# (Macroferia Regional de Talca - Poroto verde, Peru import, $/malla).
$ws.Cells.Item(140, 1).Value = 5
$ws.Cells.Item(140, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(140, 3).Value = "Maule"
$ws.Cells.Item(140, 4).Value = 44841
$ws.Cells.Item(140, 5).Value = 7
$ws.Cells.Item(140, 6).Value = 100112031
$ws.Cells.Item(140, 7).Value = "Poroto verde"
$ws.Cells.Item(140, 8).Value = "Sin especificar"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 100
$ws.Cells.Item(140, 11).Value = 33000
$ws.Cells.Item(140, 12).Value = 33000
$ws.Cells.Item(140, 13).Value = 33000
$ws.Cells.Item(140, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(140, 15).Value = "Perú"
$ws.Cells.Item(140, 16).Value = 1320
$ws.Cells.Item(140, 17).Value = 25
$ws.Cells.Item(140, 18).Value = "Hortaliza"
